$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("7").Insert()
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "20180116"
$ws.Range("B7").Value = 2148
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "20200624"
$ws.Range("B39").Value = 3230
